# Regenerate the localization handback report for the new run:
#   old GUID 87f633d4-f3dc-4621-978f-1a8809691b23 -> new GUID 3e54af86-cebf-4344-9fa1-d6578427e372
#   old commit a8545c68db732f3a4a7574ee2210eeae79ff480c -> new commit e5fa64c81e425d7da1f2f5480c85ed2d88597376
#   handoff timestamps refreshed; handback (target/handback file + handback datetime) cleared
#   since this new cycle has not been handed back yet.

$wb = $excel.ActiveWorkbook

$oldGuid = "87f633d4-f3dc-4621-978f-1a8809691b23"
$newGuid = "3e54af86-cebf-4344-9fa1-d6578427e372"
$newFileName = $newGuid + ".md"

# ---------------------------------------------------------------
# Overview sheet
# ---------------------------------------------------------------
$ws = $wb.Worksheets.Item("Overview")

$ws.Range("A2").Value = $newFileName
$ws.Range("B2").Value = "e2e\" + $newFileName
$ws.Range("G2").Value = "2016-08-28 02:59:31"

$bAddr = $ws.Range("B2").Address()
foreach ($h in $ws.Hyperlinks) {
    $addr = $h.Range.Address()
    if ($addr -eq $bAddr) {
        $h.TextToDisplay = "e2e\" + $newFileName
    }
}

# ---------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------
$ws = $wb.Worksheets.Item("zh-cn")

$ws.Range("A2").Value = $newFileName

$aAddr = $ws.Range("A2").Address()
$iAddr = $ws.Range("I2").Address()
$toDeleteHyperlinks = @()
foreach ($h in $ws.Hyperlinks) {
    $addr = $h.Range.Address()
    if ($addr -eq $aAddr) {
        $h.TextToDisplay = $newFileName
    } elseif ($addr -eq $iAddr) {
        $toDeleteHyperlinks += $h
    }
}
foreach ($h in $toDeleteHyperlinks) {
    $h.Delete()
}

$ws.Range("G2").Value = $newGuid + ".e5fa64c81e425d7da1f2f5480c85ed2d88597376.zh-cn.xlf"
$ws.Range("H2").Value = "2016-08-28 02:59:26"
$ws.Range("I2").Value = ""
$ws.Range("I2").Style = "Normal"
$ws.Range("J2").Value = ""
$ws.Range("K2").Value = "0001-01-01 00:00:00"

# The Latest Target File / Latest Handback File columns no longer hold long
# hyperlink text, so Excel shrinks them back down from the original 40-wide.
$ws.Columns.Item(9).AutoFit()
$ws.Columns.Item(10).AutoFit()

# ---------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------
$ws = $wb.Worksheets.Item("de-de")

$ws.Range("A2").Value = $newFileName

$aAddr = $ws.Range("A2").Address()
$iAddr = $ws.Range("I2").Address()
$toDeleteHyperlinks = @()
foreach ($h in $ws.Hyperlinks) {
    $addr = $h.Range.Address()
    if ($addr -eq $aAddr) {
        $h.TextToDisplay = $newFileName
    } elseif ($addr -eq $iAddr) {
        $toDeleteHyperlinks += $h
    }
}
foreach ($h in $toDeleteHyperlinks) {
    $h.Delete()
}

$ws.Range("G2").Value = $newGuid + ".e5fa64c81e425d7da1f2f5480c85ed2d88597376.de-de.xlf"
$ws.Range("H2").Value = "2016-08-28 02:59:31"
$ws.Range("I2").Value = ""
$ws.Range("I2").Style = "Normal"
$ws.Range("J2").Value = ""
$ws.Range("K2").Value = "0001-01-01 00:00:00"

$ws.Columns.Item(9).AutoFit()
$ws.Columns.Item(10).AutoFit()

Write-Output "Localization status report regenerated"
